$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7180197834968567
$ws.Range("B1").Value = 0.9313101172447205
$ws.Range("C1").Value = 0.8254098892211914
$ws.Range("D1").Value = 3.28539252281189
$ws.Range("E1").Value = 1.629402160644531
